# Mise a jour de l'application
# Add new match "N3J10" (and a couple of other late-added single-row match
# entries) data across the "Feuil1" sheet, rename the generic "CDF" header
# in HS1 to "CDF T8", and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: HS1 was the generic "CDF" label; give it its own match name ---
$ws.Range("HS1").Value = "CDF T8"

# --- New match-block data (minutes played / T-R-NR-HG result / buts / passes) ---

# Match "N3J10" -> columns EI (minutes), EJ (T/R/NR/HG), EK (buts), EL (passe D)
$ws.Range("EI2").Value = 90
$ws.Range("EJ2").Value = "T"

$ws.Range("EJ3").Value = "HG"

$ws.Range("EJ4").Value = "NR"

$ws.Range("EJ5").Value = "HG"

$ws.Range("EJ6").Value = "HG"

$ws.Range("EI7").Value = 90
$ws.Range("EJ7").Value = "T"

$ws.Range("EJ8").Value = "HG"
$ws.Range("JW8").Value = 90
$ws.Range("JX8").Value = "T"
$ws.Range("KA8").Value = 90
$ws.Range("KB8").Value = "T"

$ws.Range("EI9").Value = 90
$ws.Range("EJ9").Value = "T"

$ws.Range("EJ10").Value = "NR"

$ws.Range("EJ11").Value = "HG"
$ws.Range("KA11").Value = 90
$ws.Range("KB11").Value = "T"

$ws.Range("EJ13").Value = "HG"

$ws.Range("EI14").Value = 89
$ws.Range("EJ14").Value = "T"
$ws.Range("EK14").Value = 1
$ws.Range("EL14").Value = 1

$ws.Range("EI15").Value = 17
$ws.Range("EJ15").Value = "R"

$ws.Range("EI16").Value = 90
$ws.Range("EJ16").Value = "T"

$ws.Range("EJ17").Value = "HG"
$ws.Range("KA17").Value = 90
$ws.Range("KB17").Value = "T"

$ws.Range("EI18").Value = 73
$ws.Range("EJ18").Value = "T"
$ws.Range("EK18").Value = 1

$ws.Range("EI19").Value = 17
$ws.Range("EJ19").Value = "R"

$ws.Range("EI20").Value = 73
$ws.Range("EJ20").Value = "T"

$ws.Range("EJ21").Value = "HG"

$ws.Range("EI22").Value = 90
$ws.Range("EJ22").Value = "T"

$ws.Range("EJ23").Value = "HG"
$ws.Range("JW23").Value = 10
$ws.Range("JX23").Value = "R"

$ws.Range("EI24").Value = 90
$ws.Range("EJ24").Value = "T"
$ws.Range("EL24").Value = 1

$ws.Range("EJ25").Value = "HG"

$ws.Range("EJ26").Value = "HG"

$ws.Range("EI27").Value = 90
$ws.Range("EJ27").Value = "T"

$ws.Range("EI28").Value = 90
$ws.Range("EJ28").Value = "T"

$ws.Range("EI29").Value = 1
$ws.Range("EJ29").Value = "R"
$ws.Range("KA29").Value = 90
$ws.Range("KB29").Value = "T"

# --- View state: move the active selection the way the author left it ---
$ws.Select()
$ws.Range("EQ15").Select()
